$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Promote row 92 (last "Sandra" data row, about to become the final data row
#    after the deletions below) to the special "closing" border/format that the
#    very last data row of the table uses (copied from the current last row, 98).
$ws.Range("B98:J98").Copy()
$ws.Range("B92:J92").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Remove the two "EDER ENRIQUE RAMIREZ PEREZ" rows entirely - that worker no
#    longer appears in the statement.
$ws.Rows("16:17").Delete()

# 3) Remove the "JOHARY SARAY AREVALO MONSALVE" (3 rows), "EDIRIS GARCIA GUETO"
#    (1 row) and "LAURA TATIANA SARMIENTO ESPEJO" (2 rows) blocks - after the
#    shift from step 2 these now live at rows 91-96.
$ws.Rows("91:96").Delete()

# 4) The remaining 75 rows (16-90) are all "SANDRA PAOLA CASTRO POLO" with the
#    periods previously in descending order; rewrite them in full, in
#    chronological (ascending) order, 1611 .. 2301, each with the updated
#    arrears/salary figures.
$periods = @("1611","1612","1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112","2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212","2301")

$r = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "1046266918"
    $ws.Cells.Item($r, 4).Value = "SANDRA PAOLA CASTRO POLO"
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = 160000
    $ws.Cells.Item($r, 7).Value = 4000000
    $r++
}

# 5) Update the summary figures at the top of the statement.
$ws.Range("E11").Value = 12000000
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 75

# 6) Column D ("Nombre Trabajador") no longer needs to fit the long names that
#    were removed, so it can be narrowed to match.
$ws.Range("D:D").ColumnWidth = 28.54296875
